$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need the Text number format
# applied first, otherwise Excel auto-converts the text into a numeric value
# (e.g. losing the trailing zero in "517.08" -> 517.08 as float, or "1.00" -> 1).
$textCells = @("D5", "D6", "D10", "D12", "D13", "D16", "D19", "D21", "D23", "D24", "D26", "D28", "D31", "D32", "D34", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.954.24'
$ws.Range("E2").Value = '  -2.09%  '
$ws.Range("D3").Value = '2.464.88'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '517.08'
$ws.Range("E5").Value = '  -3.71%  '
$ws.Range("D6").Value = '131.14'
$ws.Range("E6").Value = '  -4.14%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -1.86%  '
$ws.Range("D9").Value = '2.465.12'
$ws.Range("E9").Value = '  -2.33%  '
$ws.Range("D10").Value = '0.0989'
$ws.Range("E10").Value = '  -2.31%  '
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").Value = '5.30'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '0.339'
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").Value = '2.903.50'
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").Value = '57.918.42'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '22.29'
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("E17").Value = '  -2.32%  '
$ws.Range("D18").Value = '2.474.92'
$ws.Range("E18").Value = '  -2.21%  '
$ws.Range("D19").Value = '10.71'
$ws.Range("E19").Value = '  -3.91%  '
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("D21").Value = '319.62'
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '5.71'
$ws.Range("E23").Value = '  -4.00%  '
$ws.Range("D24").Value = '64.15'
$ws.Range("E24").Value = '  -2.05%  '
$ws.Range("E25").Value = '  -2.82%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -3.72%  '
$ws.Range("D28").Value = '7.31'
$ws.Range("E28").Value = '  -2.82%  '
$ws.Range("D29").Value = '0.0₃0744'
$ws.Range("E29").Value = '  -3.53%  '
$ws.Range("E30").Value = '  -4.45%  '
$ws.Range("D31").Value = '165.22'
$ws.Range("E31").Value = '  -3.90%  '
$ws.Range("D32").Value = '6.24'
$ws.Range("E32").Value = '  -6.59%  '
$ws.Range("E33").Value = '  -1.60%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = '18.01'
$ws.Range("E36").Value = '  -2.05%  '
$ws.Range("E37").Value = '  -8.93%  '
$ws.Range("D38").Value = '3.96'
$ws.Range("E38").Value = '  -3.35%  '
$ws.Range("D39").Value = '1.46'
$ws.Range("E39").Value = '  -5.07%  '
$ws.Range("D40").Value = '0.788'
$ws.Range("E40").Value = '  -3.20%  '
$ws.Range("E41").Value = '  -4.23%  '
$ws.Range("D42").Value = '270.86'
$ws.Range("E42").Value = '  -4.98%  '
$ws.Range("D43").Value = '4.96'
$ws.Range("E43").Value = '  -2.89%  '
$ws.Range("D44").Value = '0.592'
$ws.Range("E44").Value = '  -3.04%  '
$ws.Range("D45").Value = '126.59'
$ws.Range("E45").Value = '  -3.57%  '
$ws.Range("D46").Value = '0.0906'
$ws.Range("E46").Value = '  -1.94%  '
$ws.Range("D47").Value = '0.0486'
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("E48").Value = '  -4.49%  '
$ws.Range("D49").Value = '16.85'
$ws.Range("E49").Value = '  -3.32%  '
$ws.Range("D50").Value = '1.721.15'
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("D51").Value = '0.972'
$ws.Range("E51").Value = '  -1.87%  '
